# Electricity sector calibrations
# (5) Remove CCS and H2 power plant types from RPS qualifying sources.
#
# On the "RQSD-BRQSD" (business-as-usual RPS) sheet, the following
# technologies should no longer qualify for the RPS (Boolean 1 -> 0):
#   B19  hard coal w CCS
#   B20  natural gas combined cycle w CCS
#   B21  biomass w CCS
#   B22  lignite w CCS
#   B23  small modular reactor
#   B24  hydrogen combustion turbine
#   B25  hydrogen combined cycle

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RQSD-BRQSD")
$ws.Activate()

$ws.Range("B19").Value = 0
$ws.Range("B20").Value = 0
$ws.Range("B21").Value = 0
$ws.Range("B22").Value = 0
$ws.Range("B23").Value = 0
$ws.Range("B24").Value = 0
$ws.Range("B25").Value = 0

# Leave the same range selected/highlighted, matching the author's
# on-screen state after making the edit.
$ws.Range("B19:B25").Select()
